$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(0.2439678284182306, 0.4959785522788204, 0, 0, 0, 0, 0, 0, 0.005361930294906166, 0, 0, 0, 0, 0, 0.193029490616622, 0, 0, 0.06166219839142091),
    @(0.005154639175257732, 0.03608247422680412, 0, 0, 0, 0, 0, 0, 0.01030927835051546, 0, 0, 0, 0, 0, 0.7474226804123711, 0, 0, 0.2010309278350516),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0.07407407407407407, 0, 0, 0, 0, 0, 0.6296296296296297, 0, 0, 0.2962962962962963),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0.3333333333333333, 0, 0, 0, 0, 0, 0.5, 0, 0, 0.1666666666666667),
    @(0.08133971291866028, 0, 0.01435406698564593, 0.004784688995215311, 0.06698564593301436, 0, 0, 0, 0.3205741626794258, 0, 0, 0, 0, 0.01913875598086124, 0, 0.1148325358851675, 0.08133971291866028, 0.2966507177033493),
    @(0.1106382978723404, 0, 0.01702127659574468, 0.00425531914893617, 0.0425531914893617, 0, 0, 0, 0.1574468085106383, 0, 0, 0, 0, 0.00425531914893617, 0, 0.2382978723404255, 0.08936170212765958, 0.3361702127659574),
    @(0.1081081081081081, 0, 0.01871101871101871, 0.002079002079002079, 0.04781704781704782, 0, 0, 0, 0.1143451143451144, 0, 0, 0, 0, 0.02079002079002079, 0, 0.1975051975051975, 0.103950103950104, 0.3866943866943867),
    @(0.1079545454545455, 0, 0.005681818181818182, 0, 0.03977272727272727, 0, 0, 0, 0.1420454545454546, 0, 0, 0, 0, 0.03977272727272727, 0, 0.1534090909090909, 0.125, 0.3863636363636364),
    @(0.1077119184193754, 0, 0.02613129381771829, 0.002549394518801784, 0.05162523900573614, 0, 0, 0, 0.130019120458891, 0, 0, 0, 0, 0.01720841300191205, 0, 0.2441045251752709, 0.08795411089866156, 0.3326959847036329),
    @(0, 0, 0, 0, 0, 0.1424418604651163, 0, 0, 0.0872093023255814, 0.1918604651162791, 0.563953488372093, 0, 0, 0, 0, 0, 0, 0.01453488372093023),
    @(0, 0, 0, 0, 0, 0.7598039215686274, 0, 0, 0.1568627450980392, 0.01470588235294118, 0.0392156862745098, 0, 0, 0, 0, 0, 0, 0.02941176470588235),
    @(0, 0, 0, 0, 0, 0.6122448979591837, 0, 0, 0.3469387755102041, 0, 0, 0, 0, 0, 0, 0, 0, 0.04081632653061224),
    @(0, 0, 0, 0, 0, 0.75, 0, 0, 0.25, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0.008658008658008658, 0, 0.1991341991341991, 0.0735930735930736, 0.3809523809523809, 0.05627705627705628, 0, 0, 0, 0.04761904761904762, 0, 0, 0, 0.2337662337662338),
    @(0, 0, 0, 0, 0.03734439834024896, 0, 0.1950207468879668, 0.04564315352697095, 0.4273858921161826, 0.1286307053941909, 0, 0.02074688796680498, 0, 0.04979253112033195, 0, 0, 0, 0.0954356846473029),
    @(0, 0, 0, 0, 0.01718213058419244, 0, 0.1529209621993127, 0.06872852233676977, 0.4656357388316151, 0.09450171821305842, 0, 0.02061855670103093, 0.003436426116838488, 0.0584192439862543, 0, 0, 0, 0.1185567010309278),
    @(0, 0, 0, 0, 0.0282258064516129, 0, 0.1532258064516129, 0.0564516129032258, 0.4838709677419355, 0.09274193548387097, 0, 0.01612903225806452, 0, 0.06451612903225806, 0, 0, 0, 0.1048387096774194),
    @(0, 0, 0, 0, 0.01530221882172915, 0, 0.2012241775057383, 0.07192042846212701, 0.404743687834736, 0.1124713083397093, 0, 0.02295332823259373, 0.001530221882172915, 0.05967865340474369, 0, 0, 0, 0.1101759755164499)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $data[$i][$j]
    }
}

Write-Output "Updated team specific matrix values"